$d = $word.ActiveDocument

$old = "extracted and imported in Pandas.  The relevant data were reviewed for consistency and transformed through "
$new = "extracted and imported in Jupyter notebook using Pandas.  The relevant data were reviewed for consistency and transformed through "

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
